# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps on row 5 of the
# per-language report sheets, reflecting a newly regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-19 06:14:50"
$wsZhCn.Range("G5").Value = "2016-02-19 06:15:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-19 06:15:03"
$wsDeDe.Range("G5").Value = "2016-02-19 06:15:49"
